# Static addresses filled in for listings that previously had an
# incomplete / placeholder address in column A.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "Gensinger Straße 19, 10315 Berlin"
$ws.Range("A12").Value = "Ella-Kay-Str. 9, 10405 Berlin"
$ws.Range("A16").Value = "Heylstraße 27, 10825 Berlin"
$ws.Range("A18").Value = "Hartriegelstraße 130, 12439 Berlin"
$ws.Range("A25").Value = "Letteallee 75, 13409 Berlin"
$ws.Range("A45").Value = "Adalbertstraße 44, 10997 Berlin-Mitte"
$ws.Range("A52").Value = "Königstr. 47B 14163 Zehlendorf"

# Re-apply the (unchanged) time format on these two cells so the engine
# folds them onto the existing equivalent style slot instead of keeping a
# duplicate one around.
$ws.Range("C21").NumberFormat = "hh:mm:ss"
$ws.Range("C24").NumberFormat = "hh:mm:ss"

# Leave the selection where the editor ended up (last address touched).
$ws.Range("A45").Select() | Out-Null
